$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 369
$ws.Range("A369").Value = 44810
$ws.Range("B369").Value = "KA03MN8120"
$ws.Range("C369").Value = "VENTO"
$ws.Range("D369").Value = "SUSPENSION AND PMS"
$ws.Range("E369").Value = "WORK DONE DELIVERED"
$ws.Range("F369").Value = 74840
$ws.Range("G369").Value = "CREDIT"

# Row 370
$ws.Range("A370").Value = 44779
$ws.Range("B370").Value = "KA04MM9589"
$ws.Range("C370").Value = "ECOSPORT"
$ws.Range("D370").Value = "SUSPENSION"
$ws.Range("E370").Value = "WORK DONE DELIVERED"
$ws.Range("F370").Value = 8676
$ws.Range("G370").Value = "CREDIT"

# Row 371
$ws.Range("A371").Value = 44810
$ws.Range("B371").Value = "JH09F7221"
$ws.Range("C371").Value = "SWFIT "
$ws.Range("D371").Value = "CLUTCH PROBLEM"
$ws.Range("E371").Value = "WORK DONE DELIVERED"
$ws.Range("F371").Value = 9648
$ws.Range("G371").Value = "P PAY"

# Row 372
$ws.Range("A372").Value = 44810
$ws.Range("B372").Value = "KL43B1476"
$ws.Range("C372").Value = "ALTO"
$ws.Range("D372").Value = "PMS AND CLUTCH PROBLEM"
$ws.Range("E372").Value = "WORK DONE DELIVERED"
$ws.Range("F372").Value = 9581
$ws.Range("G372").Value = "P PAY"

# Row 373 (no CASH TYPE / column G)
$ws.Range("A373").Value = 44810
$ws.Range("B373").Value = "KA25P8050"
$ws.Range("C373").Value = "NANO"
$ws.Range("D373").Value = "BATTERY CHANGE"
$ws.Range("E373").Value = "WORK DONE"
$ws.Range("F373").Value = 4956

# Row 374
$ws.Range("A374").Value = 44810
$ws.Range("B374").Value = "KA51MC4931"
$ws.Range("C374").Value = "I10 SPORTS"
$ws.Range("D374").Value = "PMS"
$ws.Range("E374").Value = "WORK DONE DELIVERED"
$ws.Range("F374").Value = 5747
$ws.Range("G374").Value = "G PAY"

# Row 375
$ws.Range("A375").Value = 44779
$ws.Range("B375").Value = "KA11B6109"
$ws.Range("C375").Value = "XYLO"
$ws.Range("D375").Value = "STARTING PROBLEM"
$ws.Range("E375").Value = "WORK DONE DELIVERED"
$ws.Range("F375").Value = 5900
$ws.Range("G375").Value = "P PAY"

# Row 376
$ws.Range("A376").Value = 44810
$ws.Range("B376").Value = "KA03MT1775"
$ws.Range("C376").Value = "XUV500"
$ws.Range("D376").Value = "WIPER BLADE CHANGE"
$ws.Range("E376").Value = "WORK DONE DELIVERED"
$ws.Range("F376").Value = 1500
$ws.Range("G376").Value = "CASH"

# Row 377
$ws.Range("A377").Value = 44810
$ws.Range("B377").Value = "KA04MJ9206"
$ws.Range("C377").Value = "CIVIC"
$ws.Range("D377").Value = "BUMPER PAINTING"
$ws.Range("E377").Value = "WORK DONE DELIVERED"
$ws.Range("F377").Value = 4000
$ws.Range("G377").Value = "P PAY"

# Update the dimension / view to match final state
$ws.Range("G377").Select()
